# 1_overview_benchmarks.xlsx - "Major change of code base (started)"
#
# Adds two new benchmark rows (DOM_GSEC / DOM_GSEC_PU, the "Domain" level)
# below the existing table, tweaks a couple of row heights / the H column
# width, reserves some extra blank (but row-formatted) rows below the new
# data, and leaves the selection parked at H24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing rows 12 & 13: row height only changes (15 -> 13.8) ---------
$ws.Rows.Item(12).RowHeight = 13.8
$ws.Rows.Item(13).RowHeight = 13.8

# --- new row 14: DOM_GSEC --------------------------------------------------
$ws.Range("A14").Value = "Domain"
$ws.Range("B14").Value = "DOM_GSEC"
$ws.Range("C14").Value = 126
$ws.Range("D14").Value = 92964
$ws.Range("E14").Value = 63
$ws.Range("F14").Value = 63
$ws.Range("H14").Value = "Prediction of gamma-secretase substrates"
$ws.Range("I14").Value = "Breimann23c"
$ws.Range("J14").Value = "1 (substrate), 0 (non-substrate)"
$ws.Rows.Item(14).RowHeight = 13.8

# --- new row 15: DOM_GSEC_PU -----------------------------------------------
$ws.Range("A15").Value = "Domain"
$ws.Range("B15").Value = "DOM_GSEC_PU"
$ws.Range("C15").Value = 694
$ws.Range("D15").Value = 494524
$ws.Range("E15").Value = 63
$ws.Range("F15").Value = 0
$ws.Range("H15").Value = "Prediction of gamma-secretase substrates (PU dataset)"
$ws.Range("I15").Value = "Breimann23c"
$ws.Range("J15").Value = "1 (substrate), 2 (unknown substrate status)"
$ws.Rows.Item(15).RowHeight = 13.8

# --- rows 16-23: reserved blank rows, formatted but with no content -------
for ($r = 16; $r -le 23; $r++) {
    $ws.Rows.Item($r).RowHeight = 13.8
}

# Touch the bottom-right corner cell so the sheet's used-range (and thus the
# <dimension> element) extends down to row 23 / out to column J, matching
# the target A1:J23 even though rows 16-23 carry no visible values.
$ws.Range("J23").NumberFormat = "General"

# --- column H width: nudge from 21.55 to 21.56 (best effort - the COM
# layer snaps widths to its internal pixel grid, so this lands on the
# nearest representable width) -----------------------------------------
$ws.Columns.Item(8).ColumnWidth = 21.56

# --- final selection --------------------------------------------------
$ws.Range("H24").Select() | Out-Null
